# Auto-generated edit script: updates crypto price/volume data
# and swaps the ImmutableX / WrappedliquidstakedEther2.0 rows (32/33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay literal text
# (matches the source data which stores prices as inline strings).
$textCells = @("D5", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values in row order.
$ws.Range("D2").Value = "24.591.63"
$ws.Range("E2").Value = "  +3.48%  "
$ws.Range("D3").Value = "1.695.95"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "316.55"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.3942"
$ws.Range("E7").Value = "  +1.52%  "
$ws.Range("D8").Value = "0.4029"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("D9").Value = "1.529"
$ws.Range("E9").Value = "  +6.29%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "53.99"
$ws.Range("E11").Value = "  +9.06%  "
$ws.Range("D12").Value = "0.08776"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").Value = "7.218"
$ws.Range("E13").Value = "  +8.17%  "
$ws.Range("D14").Value = "23.26"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "0.00001323"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "7.617"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").Value = "1.695.19"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "100.22"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "0.07055"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").Value = "6.856"
$ws.Range("E21").Value = "  +2.98%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "14.08"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Value = "24.576.97"
$ws.Range("E24").Value = "  +3.48%  "
$ws.Range("D25").Value = "3.018"
$ws.Range("E25").Value = "  +7.58%  "
$ws.Range("D26").Value = "2.308"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "22.43"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("D28").Value = "159.88"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").Value = "5.227"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "134.76"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").Value = "7.505"
$ws.Range("E31").Value = "  +16.38%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.111"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.881.79"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").Value = "7.363"
$ws.Range("E34").Value = "  +12.29%  "
$ws.Range("D35").Value = "0.08532"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "11.42"
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("D37").Value = "1.958"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "0.2741"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("D39").Value = "14.57"
$ws.Range("E39").Value = "  +0.72%  "
$ws.Range("D40").Value = "0.02764"
$ws.Range("E40").Value = "  +9.91%  "
$ws.Range("D41").Value = "0.09076"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "1.466"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "0.7693"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").Value = "0.7182"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").Value = "15.51"
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("D46").Value = "2.527"
$ws.Range("E46").Value = "  +5.72%  "
$ws.Range("D47").Value = "4.215"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "1.355"
$ws.Range("E48").Value = "  +12.65%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "141.27"
$ws.Range("E50").Value = "  +2.25%  "
$ws.Range("D51").Value = "0.08025"
$ws.Range("E51").Value = "  +3.45%  "
